$wb = $excel.ActiveWorkbook

# Update "OFF" sheet (Home row, r=2): simulated/logged 2021 conference championship stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 419
$wsOff.Range("C2").Value = 294
$wsOff.Range("D2").Value = 91
$wsOff.Range("E2").Value = 41
$wsOff.Range("F2").Value = 7
$wsOff.Range("G2").Value = 4

# Update "DEF" sheet (Home row, r=2): simulated/logged 2021 conference championship stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 401
$wsDef.Range("C2").Value = 282
$wsDef.Range("D2").Value = 84
$wsDef.Range("E2").Value = 29
$wsDef.Range("F2").Value = 8
$wsDef.Range("G2").Value = 5
